$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.0041407867494824
$ws.Range("D2").Value = 0.0020703933747412
$ws.Range("E2").Value = 0.0144927536231884
$ws.Range("F2").Value = 0.0020703933747412
$ws.Range("H2").Value = 0.0020703933747412
$ws.Range("J2").Value = 0.0062111801242236
$ws.Range("L2").Value = 0.0041407867494824
$ws.Range("M2").Value = 0.975155279503106
$ws.Range("N2").Value = 0.0289855072463768
$ws.Range("O2").Value = 0.958592132505176
$ws.Range("P2").Value = 0.954451345755694
$ws.Range("Q2").Value = 0.0020703933747412
$ws.Range("S2").Value = 0.989648033126294
$ws.Range("T2").Value = 0.0020703933747412
$ws.Range("U2").Value = 0.995859213250518
$ws.Range("W2").Value = 0.884057971014493
$ws.Range("X2").Value = 0.062111801242236
$ws.Range("B3").Value = 0.987577639751553
$ws.Range("C3").Value = 0.981366459627329
$ws.Range("D3").Value = 0.985507246376812
$ws.Range("E3").Value = 0.0082815734989648
$ws.Range("H3").Value = 0.0248447204968944
$ws.Range("I3").Value = 0.0020703933747412
$ws.Range("J3").Value = 0.0020703933747412
$ws.Range("L3").Value = 0.98343685300207
$ws.Range("M3").Value = 0.0041407867494824
$ws.Range("N3").Value = 0.0062111801242236
$ws.Range("P3").Value = 0.010351966873706
$ws.Range("Q3").Value = 0.946169772256729
$ws.Range("R3").Value = 0.010351966873706
$ws.Range("T3").Value = 0.989648033126294
$ws.Range("U3").Value = 0.0020703933747412
$ws.Range("V3").Value = 0.985507246376812
$ws.Range("W3").Value = 0.0993788819875776
$ws.Range("X3").Value = 0.921325051759834
$ws.Range("B4").Value = 0.0041407867494824
$ws.Range("C4").Value = 0.0082815734989648
$ws.Range("D4").Value = 0.0082815734989648
$ws.Range("E4").Value = 0.0082815734989648
$ws.Range("F4").Value = 0.997929606625259
$ws.Range("H4").Value = 0.010351966873706
$ws.Range("J4").Value = 0.991718426501035
$ws.Range("L4").Value = 0.010351966873706
$ws.Range("M4").Value = 0.0186335403726708
$ws.Range("N4").Value = 0.954451345755694
$ws.Range("O4").Value = 0.010351966873706
$ws.Range("P4").Value = 0.031055900621118
$ws.Range("R4").Value = 0.0041407867494824
$ws.Range("S4").Value = 0.0020703933747412
$ws.Range("T4").Value = 0.0020703933747412
$ws.Range("U4").Value = 0.0020703933747412
$ws.Range("V4").Value = 0.0124223602484472
$ws.Range("W4").Value = 0.0144927536231884
$ws.Range("X4").Value = 0.0020703933747412
$ws.Range("B5").Value = 0.0062111801242236
$ws.Range("C5").Value = 0.0062111801242236
$ws.Range("D5").Value = 0.0041407867494824
$ws.Range("E5").Value = 0.968944099378882
$ws.Range("H5").Value = 0.962732919254658
$ws.Range("I5").Value = 0.997929606625259
$ws.Range("L5").Value = 0.0020703933747412
$ws.Range("M5").Value = 0.0020703933747412
$ws.Range("N5").Value = 0.010351966873706
$ws.Range("O5").Value = 0.0289855072463768
$ws.Range("P5").Value = 0.0020703933747412
$ws.Range("Q5").Value = 0.05175983436853
$ws.Range("R5").Value = 0.985507246376812
$ws.Range("S5").Value = 0.0082815734989648
$ws.Range("T5").Value = 0.0062111801242236
$ws.Range("V5").Value = 0.0020703933747412
$ws.Range("W5").Value = 0.0020703933747412
$ws.Range("X5").Value = 0.0144927536231884
